# Weekly data refresh: a new price-report row for "Zapallo italiano" was
# added to the Macroferia Regional de Talca sheet. The new observation is
# inserted as row 129 (pushing the existing rows 129-192 down to 130-193),
# matching how the source data is prepended with each week's latest entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 129, shifting rows 129:192 down to 130:193.
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row with this week's reported values.
$ws.Range("A129").Value = 5
$ws.Range("B129").Value = "Macroferia Regional de Talca"
$ws.Range("C129").Value = "Maule"
$ws.Range("D129").Value = 44455
$ws.Range("E129").Value = 7
$ws.Range("F129").Value = 100112032
$ws.Range("G129").Value = "Zapallo italiano"
$ws.Range("H129").Value = "Sin especificar"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 200
$ws.Range("K129").Value = 16000
$ws.Range("L129").Value = 16000
$ws.Range("M129").Value = 16000
$ws.Range("N129").Value = "`$/caja 50 unidades"
$ws.Range("O129").Value = "Región de Arica y Parinacota"
$ws.Range("P129").Value = 320
$ws.Range("Q129").Value = 50
$ws.Range("R129").Value = "Hortaliza"
